$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("GDPbES")

# ---------------------------------------------------------------------------
# 1. Update the guaranteed-dispatch percentages in column B (rows 3-10).
#    Every other year-column (C:AK) already holds a "=$B<row>" shared
#    formula, so writing the literal value into B recalculates the rest of
#    the row automatically.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 0.2    # natural gas nonpeaker
$ws.Range("B4").Value = 0.5    # nuclear
$ws.Range("B5").Value = 0.5    # hydro
$ws.Range("B6").Value = 0.8    # onshore wind
$ws.Range("B7").Value = 0.8    # solar PV
$ws.Range("B8").Value = 0.8    # solar thermal
$ws.Range("B9").Value = 0.5    # biomass
$ws.Range("B10").Value = 0.8   # geothermal

# ---------------------------------------------------------------------------
# 2. Add three new fuel-source rows (15-17) that mirror existing rows via
#    formula, same pattern as the existing sheet (row formula references a
#    sibling row, then fills right with a shared formula).
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"
$ws.Range("C15:AK15").Formula = "=C11"

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"
$ws.Range("C16:AK16").Formula = "=C11"

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"
$ws.Range("C17:AK17").Formula = "=C9"

# ---------------------------------------------------------------------------
# 3. Give row 1 a label in column A describing the units, bold + wrapped,
#    with a taller row to fit the wrapped text.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Guaranteed Dispatch Fraction (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

# ---------------------------------------------------------------------------
# 4. Switch the active tab from GDPbES to About.
# ---------------------------------------------------------------------------
[void]$ws.Range("A1").Select()
[void]$wsAbout.Activate()
[void]$wsAbout.Range("A1").Select()
